# Fill in the weekly fantasy-football scores on the "Scores" sheet.
# Rows 2-4 already had week1-week8 (cols B-I); this adds week9-week11
# (cols J-L). Rows 5-13 had no weekly scores at all yet; this fills in
# the full week1-week11 range (cols B-L) for every team.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

$row2 = @(122.3, 123.9, 132.6, 108.2, 68.8, 63, 87.1, 54.7, 101.5, 65.2, 92)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value2 = $row2[$i]
}

$row3 = @(120, 119.5, 110, 74.3, 95.7, 117.4, 127.1, 103.1, 52.8, 155.19999999999999, 122.6)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 2 + $i).Value2 = $row3[$i]
}

$row4 = @(70.8, 73.7, 122, 79.3, 111, 131.69999999999999, 70.3, 117.4, 80.900000000000006, 63.4, 107.8)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 2 + $i).Value2 = $row4[$i]
}

$row5 = @(109.2, 98.2, 89, 92.8, 124, 126.9, 90.8, 103.8, 97.2, 90.3, 118)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 2 + $i).Value2 = $row5[$i]
}

$row6 = @(81.8, 96.6, 91.8, 98.4, 123.8, 97.7, 94.6, 83.1, 126.1, 75, 90)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, 2 + $i).Value2 = $row6[$i]
}

$row7 = @(82.8, 90.4, 125.7, 112.9, 81, 91.7, 95.4, 123.3, 115.9, 58.2, 115.3)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, 2 + $i).Value2 = $row7[$i]
}

$row8 = @(82.2, 75.7, 68.599999999999994, 95.1, 81.900000000000006, 51.3, 95.3, 80.099999999999994, 102.1, 102.3, 51.8)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, 2 + $i).Value2 = $row8[$i]
}

$row9 = @(88.1, 118.6, 74.5, 125.1, 100.7, 91.6, 82.2, 122.1, 111.1, 74.900000000000006, 106.9)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, 2 + $i).Value2 = $row9[$i]
}

$row10 = @(111.4, 93.4, 116, 86.6, 110.5, 87.7, 116.3, 96.2, 106.8, 105.6, 81.7)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, 2 + $i).Value2 = $row10[$i]
}

$row11 = @(64.2, 89, 71.099999999999994, 104.6, 111.6, 99.5, 84.2, 126.5, 81.2, 106.3, 94)
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, 2 + $i).Value2 = $row11[$i]
}

$row12 = @(72.3, 109.8, 134, 116.2, 88.4, 98.9, 97.4, 97.3, 93.3, 65.5, 107)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, 2 + $i).Value2 = $row12[$i]
}

$row13 = @(113.8, 121.1, 104.1, 96.4, 116.4, 108.4, 105.2, 107.8, 90.2, 116.2, 114.5)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $ws.Cells.Item(13, 2 + $i).Value2 = $row13[$i]
}

# Match the author's final active-cell selection on the Scores sheet.
$ws.Range("I12").Select() | Out-Null
